# [Fonds de solidarite] Add 2022-05-19 data
# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# for the rows whose underlying aggregated figures changed with the new
# 2022-05-19 data extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 10;  C = 345539;  D = 64165;  E = 1817720252 },
    @{ Row = 21;  C = 175240;  D = 38057;  E = 316813398 },
    @{ Row = 78;  C = 178441;  D = 34685;  E = 892526355 },
    @{ Row = 121; C = 1306164; D = 220385; E = 2274669692 },
    @{ Row = 129; C = 633440;  D = 104966; E = 3428621480 },
    @{ Row = 132; C = 585716;  D = 90780;  E = 3463679198 },
    @{ Row = 144; C = 25079;   D = 6170;   E = 92448942 },
    @{ Row = 154; C = 18455;   D = 3296;   E = 72768683 },
    @{ Row = 156; C = 12402;   D = 2144;   E = 40299139 },
    @{ Row = 194; C = 18380;   D = 2990;   E = 71349778 },
    @{ Row = 229; C = 612546;  D = 121244; E = 1040763491 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
